$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename Sheet1 -> "TEST 20-12"
$ws.Name = "TEST 20-12"

# Highlight flagged question cells with a solid red fill (tag-style marking
# to help filter questions later)
$ws.Range("B5").Interior.Color = 255
$ws.Range("B7").Interior.Color = 255
$ws.Range("B9").Interior.Color = 255

# Remove the now-unused "answer" column C (was only ever a stray header +
# one empty styled cell) and tidy the selection/scroll position
$ws.Columns.Item(3).Delete()
$ws.Range("B18").Select()
